# Added Buy Deals for Inter and Euro resellers
#
# 1. Rename existing "NymgoNormalEuroUser" sheet to "NymgoEuroNormalUser".
# 2. Duplicate it into a new sheet "NymgoEuroReseller" placed right after it.
# 3. Populate the new sheet's account-details column (B) with the reseller
#    deal-account data (mirrors the "dealtester" -> "dealReseller" pattern
#    already used by the sibling "Deal Tester" sheet).
# 4. Leave the new sheet active/selected, matching the tab state Excel
#    leaves behind after such an edit.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("NymgoNormalEuroUser")

# Rename the template sheet first (the new copy is derived from it, so the
# copy's own name is set afterwards).
$template.Name = "NymgoEuroNormalUser"

# Duplicate the sheet immediately after itself.
$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item($template.Index + 1)
$newSheet.Name = "NymgoEuroReseller"

# Username
$newSheet.Range("B1").Value = "dealReseller"
# Password stays "password" (unchanged, already copied)
# FullName
$newSheet.Range("B3").Value = "Deal Reseller Account"
# Email stays "dealtester@mail.ru" (unchanged, already copied)
# Mobile
$newSheet.Range("B5").Value = "111111111"
# Phone (was a formula referencing TestAccountsInfo!H9 - now a literal)
$newSheet.Range("B6").Value = "111111111"
# CountryOfResidence (was a formula referencing TestAccountsInfo!I9 - now a literal)
$newSheet.Range("B7").Value = "Finland"
# VAT
$newSheet.Range("B16").Value = "24"
# BusinessName
$newSheet.Range("B19").Value = "Business name"

# Leave the freshly-created reseller sheet selected/active, with the cursor
# on B17 (matches the recorded selection after populating the sheet).
$newSheet.Activate() | Out-Null
$newSheet.Range("B17").Select() | Out-Null
